$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new column before column A, shifting every existing
# column (and the data validation that targets the last column) one to
# the right. Excel keeps per-cell formatting intact across the insert.
$ws.Columns("A").Insert()

# Give the new column A the same look as its header/data neighbours by
# copying formats only, then fill in the new values.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"

$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 3

# Column width for the new column A (stored width 23).
$ws.Columns("A").ColumnWidth = 22.16666666666667

# Upper-case every existing header label (now living in B1:X1); the
# right-most status column (now Y1) is left exactly as it was.
$ws.Range("B1").Value = "REGION"
$ws.Range("C1").Value = "DIVISION"
$ws.Range("D1").Value = "SCHOOL ID"
$ws.Range("E1").Value = "SCHOOL NAME"
$ws.Range("F1").Value = "MUNICIPALITY"
$ws.Range("G1").Value = "LEG DISTRICT"
$ws.Range("H1").Value = "NO. OF SITES"
$ws.Range("I1").Value = "SCOPE OF WORK"
$ws.Range("J1").Value = "TOTAL ALLOCATION"
$ws.Range("K1").Value = "CONTRACT AMOUNT"
$ws.Range("L1").Value = "STATUS"
$ws.Range("M1").Value = "PERCENTAGE OF COMPLETION"
$ws.Range("N1").Value = " TARGET COMPLETION DATE "
$ws.Range("O1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("P1").Value = "PROJECT ID"
$ws.Range("Q1").Value = "CONTRACT ID"
$ws.Range("R1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("S1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("T1").Value = "BID OPENING"
$ws.Range("U1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("V1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("W1").Value = "NAME OF CONTRACTOR"
$ws.Range("X1").Value = "OTHER REMARKS"
